$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 18
$ws.Range("H18").Value = 299.66666
$ws.Range("I18").Value = 299.66666
$ws.Range("K18").Value = 299.66666
$ws.Range("M18").Value = -15.66665999999998
# Row 19
$ws.Range("H19").Value = 1268.65
$ws.Range("I19").Value = 837
$ws.Range("J19").Value = 1556.4166
$ws.Range("K19").Value = 837
$ws.Range("L19").Value = 1556.4166
$ws.Range("M19").Value = -662
$ws.Range("N19").Value = -1906.4166
# Row 41
$ws.Range("H41").Value = 1179.5714
$ws.Range("I41").Value = 1467
$ws.Range("K41").Value = 1467
$ws.Range("M41").Value = -1027
# Row 43
$ws.Range("H43").Value = 8718.143
$ws.Range("J43").Value = 1749.3636
$ws.Range("L43").Value = 1749.3636
$ws.Range("N43").Value = -1887.3636
# Row 45
$ws.Range("J45").Value = 18.5
$ws.Range("L45").Value = 55.5
$ws.Range("N45").Value = -439.5
# Row 98
$ws.Range("H98").Value = 861.8095
$ws.Range("I98").Value = 829.2941
$ws.Range("K98").Value = 829.2941
$ws.Range("M98").Value = 668.7059
# Row 107
$ws.Range("H107").Value = 2075.7334
$ws.Range("I107").Value = 648.8182
$ws.Range("K107").Value = 648.8182
$ws.Range("M107").Value = 1271.1818
# Row 122
$ws.Range("H122").Value = 861.8095
$ws.Range("I122").Value = 829.2941
$ws.Range("K122").Value = 2487.8823
$ws.Range("M122").Value = -37.88229999999976
# Row 137
$ws.Range("H137").Value = 1071.6207
$ws.Range("I137").Value = 831.36365
$ws.Range("K137").Value = 2494.09095
$ws.Range("M137").Value = 55.90905000000021
# Row 138
$ws.Range("H138").Value = 4373.6
$ws.Range("J138").Value = 5647.826
$ws.Range("L138").Value = 16943.478
$ws.Range("N138").Value = -27223.478

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 2916.923
$ws.Range("I2").Value = 3874.2856
$ws.Range("J2").Value = 1800
$ws.Range("K2").Value = 3874.2856
$ws.Range("L2").Value = 1800
$ws.Range("M2").Value = -3761.2856
$ws.Range("N2").Value = -2026
# Row 32
$ws.Range("H32").Value = 3578.6177
$ws.Range("I32").Value = 3183.5312
$ws.Range("K32").Value = 3183.5312
$ws.Range("M32").Value = -2896.5312
# Row 45
$ws.Range("H45").Value = 23337.916
$ws.Range("I45").Value = 23220.143
$ws.Range("K45").Value = 23220.143
$ws.Range("M45").Value = -22843.143
# Row 61
$ws.Range("H61").Value = 2689.5483
$ws.Range("I61").Value = 2614.52
$ws.Range("J61").Value = 3002.1667
$ws.Range("K61").Value = 2614.52
$ws.Range("L61").Value = 3002.1667
$ws.Range("M61").Value = -2402.52
$ws.Range("N61").Value = -3426.1667
# Row 102
$ws.Range("H102").Value = 3392.5715
$ws.Range("I102").Value = 2949.6
$ws.Range("K102").Value = 2949.6
$ws.Range("M102").Value = -1327.6
# Row 116
$ws.Range("H116").Value = 2916.923
$ws.Range("I116").Value = 3874.2856
$ws.Range("J116").Value = 1800
$ws.Range("K116").Value = 3874.2856
$ws.Range("L116").Value = 1800
$ws.Range("M116").Value = -1580.2856
$ws.Range("N116").Value = -6388
# Row 132
$ws.Range("H132").Value = 3560.2727
$ws.Range("I132").Value = 3490.9473
$ws.Range("J132").Value = 3999.3333
$ws.Range("K132").Value = 10472.8419
$ws.Range("L132").Value = 11997.9999
$ws.Range("M132").Value = -7942.841899999999
$ws.Range("N132").Value = -17057.9999
# Row 136
$ws.Range("H136").Value = 2689.5483
$ws.Range("I136").Value = 2614.52
$ws.Range("J136").Value = 3002.1667
$ws.Range("K136").Value = 7843.559999999999
$ws.Range("L136").Value = 9006.500100000001
$ws.Range("M136").Value = -5293.559999999999
$ws.Range("N136").Value = -14106.5001
# Row 139
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 2916.923
$ws.Range("I3").Value = 3874.2856
$ws.Range("J3").Value = 1800
$ws.Range("K3").Value = 3874.2856
$ws.Range("L3").Value = 1800
$ws.Range("M3").Value = -3760.2856
$ws.Range("N3").Value = -2028
# Row 107
$ws.Range("H107").Value = 4414.0713
$ws.Range("I107").Value = 3533.111
$ws.Range("J107").Value = 5999.8
$ws.Range("K107").Value = 3533.111
$ws.Range("L107").Value = 5999.8
$ws.Range("M107").Value = -1613.111
$ws.Range("N107").Value = -9839.799999999999

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 32
$ws.Range("H32").Value = 3411.3333
$ws.Range("I32").Value = 3411.3333
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 3411.3333
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -3095.3333
$ws.Range("N32").ClearContents()
# Row 107
$ws.Range("H107").Value = 1028.6666
$ws.Range("I107").Value = 984.4
$ws.Range("J107").Value = 1250
$ws.Range("K107").Value = 984.4
$ws.Range("L107").Value = 1250
$ws.Range("M107").Value = 935.6
$ws.Range("N107").Value = -5090

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 32
$ws.Range("H32").Value = 1995.5
$ws.Range("J32").Value = 1990
$ws.Range("L32").Value = 5970
$ws.Range("N32").Value = -6536
# Row 92
$ws.Range("H92").Value = 652.1667
$ws.Range("I92").Value = 271.83334
$ws.Range("J92").Value = 842.3333
$ws.Range("K92").Value = 815.5000200000001
$ws.Range("L92").Value = 2526.9999
$ws.Range("M92").Value = 432.4999799999999
$ws.Range("N92").Value = -5022.9999

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 12
$ws.Range("H12").Value = 63335.332
$ws.Range("I12").Value = 63335.332
$ws.Range("K12").Value = 63335.332
$ws.Range("M12").Value = -63195.332
# Row 80
$ws.Range("H80").Value = 2850.7144
$ws.Range("I80").Value = 2816.6667
$ws.Range("J80").Value = 2876.25
$ws.Range("K80").Value = 2816.6667
$ws.Range("L80").Value = 2876.25
$ws.Range("M80").Value = -1818.6667
$ws.Range("N80").Value = -4872.25
# Row 83
$ws.Range("H83").Value = 2850.7144
$ws.Range("I83").Value = 2816.6667
$ws.Range("J83").Value = 2876.25
$ws.Range("K83").Value = 14083.3335
$ws.Range("L83").Value = 14381.25
$ws.Range("M83").Value = -9091.333500000001
$ws.Range("N83").Value = -24365.25
# Row 126
$ws.Range("H126").Value = 25714.857
$ws.Range("J126").Value = 3157.1428
$ws.Range("L126").Value = 9471.428400000001
$ws.Range("N126").Value = -14411.4284
# Row 132
$ws.Range("H132").Value = 259078.84
$ws.Range("I132").Value = 272915.53
$ws.Range("J132").Value = 3100
$ws.Range("K132").Value = 818746.5900000001
$ws.Range("L132").Value = 9300
$ws.Range("M132").Value = -816216.5900000001
$ws.Range("N132").Value = -14360

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 6
$ws.Range("H6").Value = 48500
$ws.Range("J6").Value = 48500
$ws.Range("L6").Value = 48500
$ws.Range("N6").Value = -48724
# Row 7
$ws.Range("H7").Value = 24214.715
$ws.Range("I7").Value = 26750.5
$ws.Range("K7").Value = 26750.5
$ws.Range("M7").Value = -26638.5
# Row 16
$ws.Range("H16").Value = 3964.9355
$ws.Range("I16").Value = 1581.5333
$ws.Range("K16").Value = 1581.5333
$ws.Range("M16").Value = -1411.5333
# Row 61
$ws.Range("H61").Value = 4360.1113
$ws.Range("I61").Value = 4544.727
$ws.Range("J61").Value = 4070
$ws.Range("K61").Value = 4544.727
$ws.Range("L61").Value = 4070
$ws.Range("M61").Value = -4342.727
$ws.Range("N61").Value = -4474
# Row 68
$ws.Range("H68").Value = 5752.625
$ws.Range("I68").Value = 3426.9412
$ws.Range("J68").Value = 11400.714
$ws.Range("K68").Value = 3426.9412
$ws.Range("L68").Value = 11400.714
$ws.Range("M68").Value = -2677.9412
$ws.Range("N68").Value = -12898.714
# Row 71
$ws.Range("H71").Value = 5752.625
$ws.Range("I71").Value = 3426.9412
$ws.Range("J71").Value = 11400.714
$ws.Range("K71").Value = 17134.706
$ws.Range("L71").Value = 57003.57
$ws.Range("M71").Value = -13390.706
$ws.Range("N71").Value = -64491.57
# Row 96
$ws.Range("H96").Value = 49850
$ws.Range("J96").Value = 49850
$ws.Range("L96").Value = 49850
$ws.Range("N96").Value = -55342
# Row 104
$ws.Range("H104").Value = 50001
$ws.Range("J104").Value = 50001
$ws.Range("L104").Value = 50001
$ws.Range("N104").Value = -56989
# Row 113
$ws.Range("H113").Value = 4360.1113
$ws.Range("I113").Value = 4544.727
$ws.Range("J113").Value = 4070
$ws.Range("K113").Value = 4544.727
$ws.Range("L113").Value = 4070
$ws.Range("M113").Value = -2374.727
$ws.Range("N113").Value = -8410
# Row 122
$ws.Range("H122").Value = 483296.38
$ws.Range("I122").Value = 592539.7
$ws.Range("K122").Value = 1777619.1
$ws.Range("M122").Value = -1775169.1
# Row 126
$ws.Range("H126").Value = 24214.715
$ws.Range("I126").Value = 26750.5
$ws.Range("K126").Value = 80251.5
$ws.Range("M126").Value = -77781.5
# Row 136
$ws.Range("H136").Value = 4641.533
$ws.Range("I136").Value = 4630.9546
$ws.Range("J136").Value = 4670.625
$ws.Range("K136").Value = 13892.8638
$ws.Range("L136").Value = 14011.875
$ws.Range("M136").Value = -11342.8638
$ws.Range("N136").Value = -19111.875
# Row 141
$ws.Range("H141").Value = 84113.336
$ws.Range("J141").Value = 84113.336
$ws.Range("L141").Value = 84113.336
$ws.Range("N141").Value = -94473.336

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 41
$ws.Range("H41").Value = 20816.428
$ws.Range("I41").Value = 21500
$ws.Range("J41").Value = 20543
$ws.Range("K41").Value = 20543
$ws.Range("L41").Value = 20543
$ws.Range("M41").Value = -21110
$ws.Range("N41").Value = -21323
# Row 126
$ws.Range("H126").Value = 7149.577
$ws.Range("I126").Value = 7585.409
$ws.Range("K126").Value = 22756.227
$ws.Range("M126").Value = -20286.227
# Row 132
$ws.Range("H132").Value = 3429.5264
$ws.Range("I132").Value = 2823.1875
$ws.Range("K132").Value = 8469.5625
$ws.Range("M132").Value = -5939.5625
# Row 136
$ws.Range("H136").Value = 1541.9032
$ws.Range("I136").Value = 1251.8889
$ws.Range("K136").Value = 3755.6667
$ws.Range("M136").Value = -1205.6667
